$wb = $excel.ActiveWorkbook

$wsYearly = $wb.Worksheets.Item("Yearly")
$wsAllTime = $wb.Worksheets.Item("All Time")

# The October (row 12) "Taxable Account" dividend amount was corrected.
$wsYearly.Range("D12").Value = 35.24

# The "All Time" summary sheet keeps a manually entered copy of the 2016
# "Taxable Account" yearly total (row 7) - update it to match.
$wsAllTime.Range("F7").Value = 608.14

# Recalculate all dependent totals (Grand Total / SUM formulas) on both
# sheets so the cached formula results stay in sync with the new inputs.
$excel.CalculateFullRebuild()

# Restore the on-screen selection state that was captured for each sheet.
$wsYearly.Activate()
$wsYearly.Range("I13").Select()

$wsAllTime.Activate()
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
$wsAllTime.Range("J52").Select()
